$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 80 (shifts existing rows 80:107 down to 81:108)
$ws.Rows.Item(80).Insert()

# Populate the newly inserted row 80 with this week's price report entry
$ws.Cells.Item(80, 1).Value = 9
$ws.Cells.Item(80, 2).Value = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(80, 3).Value = "Metropolitana"
$ws.Cells.Item(80, 4).Value = 45135
$ws.Cells.Item(80, 5).Value = 13
$ws.Cells.Item(80, 6).Value = 100112035
$ws.Cells.Item(80, 7).Value = "Bruselas (repollito)"
$ws.Cells.Item(80, 8).Value = "Sin especificar"
$ws.Cells.Item(80, 9).Value = "Primera"
$ws.Cells.Item(80, 10).Value = 52
$ws.Cells.Item(80, 11).Value = 17000
$ws.Cells.Item(80, 12).Value = 18000
$ws.Cells.Item(80, 13).Value = 17500
$ws.Cells.Item(80, 14).Value = "$/malla 15 kilos"
$ws.Cells.Item(80, 15).Value = "Provincia de Quillota"
$ws.Cells.Item(80, 16).Value = 1167
$ws.Cells.Item(80, 17).Value = 15
$ws.Cells.Item(80, 18).Value = "Hortaliza"
